$d = $word.ActiveDocument

# --- 1) "InstantPay" bullet -> rewritten text/formatting, numId 1 -> 2 ---
$xmlInstantPay = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
      '<w:color w:val="252525"/>' +
      '<w:sz w:val="21"/>' +
      '<w:szCs w:val="21"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
      '<w:color w:val="252525"/>' +
      '<w:sz w:val="21"/>' +
      '<w:szCs w:val="21"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t>' +
  '</w:r>' +
'</w:p>'

# --- 2) "Will have more services added later" bullet -> SmartRewards text/formatting ---
$xmlSmartRewards = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
      '<w:color w:val="252525"/>' +
      '<w:sz w:val="21"/>' +
      '<w:szCs w:val="21"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
      '<w:color w:val="252525"/>' +
      '<w:sz w:val="21"/>' +
      '<w:szCs w:val="21"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t>' +
  '</w:r>' +
'</w:p>'

# Locate the two bullet paragraphs by their current (pre-edit) text so the
# script is resilient to being re-run against the same starting document.
$targetInstantPay = $null
$targetServices = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("InstantPay") -and $targetInstantPay -eq $null) {
        $targetInstantPay = $i
    }
    if ($t.StartsWith("Will have more services") -and $targetServices -eq $null) {
        $targetServices = $i
    }
}

$pInstantPay = $d.Paragraphs.Item($targetInstantPay)
$pInstantPay.Range.InsertXML($xmlInstantPay)

$pServices = $d.Paragraphs.Item($targetServices)
$pServices.Range.InsertXML($xmlSmartRewards)
